$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns E and F (Effort, SLOC)
$ws.Range("E2").Value = "Effort"
$ws.Range("F2").Value = "SLOC"

# Column E (Effort) - fill in missing "NA" values for rows that previously had none
$ws.Range("E5").Value = "NA"

# Column F (SLOC) - new data column
$ws.Range("F5").Value = "NA"
$ws.Range("F6").Value = "NA"
$ws.Range("F7").Value = "NA"
$ws.Range("F8").Value = "NA"
$ws.Range("F9").Value = "NA"
$ws.Range("F10").Value = 5400
$ws.Range("F11").Value = "NA"
$ws.Range("F12").Value = 1393
$ws.Range("F13").Value = 2042
$ws.Range("F14").Value = 1775
$ws.Range("F15").Value = 4965
$ws.Range("F16").Value = 3795
$ws.Range("F17").Value = 3000
$ws.Range("F18").Value = "NA"
$ws.Range("F19").Value = 7776
$ws.Range("F20").Value = 4095
$ws.Range("F21").Value = "NA"
$ws.Range("F22").Value = 3507
$ws.Range("F23").Value = 5330
$ws.Range("F24").Value = 3600
$ws.Range("F25").Value = "NA"

# Fix a data point: team no. for row 11 changed from 2 to 1
$ws.Range("A11").Value = 1

# Column width adjustments (D and E narrower now that F holds SLOC data)
$ws.Columns.Item(4).ColumnWidth = 48
$ws.Columns.Item(5).ColumnWidth = 22

# Update the selected cell (was C17, now C16)
$ws.Range("C16").Select()
